# Daily attendance processing - rotate the "Recorded By" (column G) author
# list so the first-listed recorder moves to the end of the comma-separated
# list. Single-author cells are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count + $usedRange.Row - 1

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $val = $cell.Value2

    if ($val -eq $null) { continue }

    $parts = $val -split ", "

    if ($parts.Length -gt 1) {
        $rotated = $parts[1..($parts.Length - 1)] + @($parts[0])
        $cell.Value = $rotated -join ", "
    }
}
